# Updated test data as per new implementation:
# - Move the "MX 4000 / CPU 801" data row (old row 8) out to a brand new
#   "Test data" worksheet (as its row 1), leaving only the
#   "Pro32xD / PFI" row behind on "Add Panels" (now row 8, shifted up
#   from row 9).
# - Relabel the M/N values on the remaining row (now row 8) from
#   "Battery Alarm (A)" / "Battery Standby (A)" to the new
#   "Alarm Current(A)" / "Standby Current(A)" labels.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet right after "Add Panels".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Test data"

# Move (cut) the old row 8 (MX 4000 / CPU 801 / FIM data) into the new sheet,
# then delete the now-empty row so row 9 shifts up to become row 8.
$ws1.Range("A8:N8").Cut($ws2.Range("A1:N1"))
$ws1.Rows(8).Delete()

# Update the relabeled header values on the row that is now row 8.
$ws1.Range("M8").Value = "Alarm Current(A)"
$ws1.Range("N8").Value = "Standby Current(A)"

# Fix up the selections: "Test data" keeps its whole first row selected,
# "Add Panels" selects M8:N8 and stays the active/visible tab.
$ws2.Rows(1).Select()
$ws1.Range("M8:N8").Select()
$ws1.Activate()
